$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1066.5
$ws.Range("I18").Value = 1066.5
$ws.Range("K18").Value = 1066.5
$ws.Range("M18").Value = -782.5
$ws.Range("H43").Value = 1858.4375
$ws.Range("I43").Value = 969.8889
$ws.Range("J43").Value = 3000.8572
$ws.Range("K43").Value = 969.8889
$ws.Range("L43").Value = 3000.8572
$ws.Range("M43").Value = -900.8889
$ws.Range("N43").Value = -3138.8572
$ws.Range("H76").Value = 3344.7778
$ws.Range("I76").Value = 3283.8333
$ws.Range("J76").Value = 3466.6667
$ws.Range("K76").Value = 3283.8333
$ws.Range("L76").Value = 3466.6667
$ws.Range("M76").Value = -2968.8333
$ws.Range("N76").Value = -4096.6667
$ws.Range("H79").Value = 3344.7778
$ws.Range("I79").Value = 3283.8333
$ws.Range("J79").Value = 3466.6667
$ws.Range("K79").Value = 3283.8333
$ws.Range("L79").Value = 3466.6667
$ws.Range("M79").Value = -2191.8333
$ws.Range("N79").Value = -5650.6667
$ws.Range("H132").Value = 1091016.1
$ws.Range("I132").Value = 2013.4324
$ws.Range("K132").Value = 6040.2972
$ws.Range("M132").Value = -3510.2972
$ws.Range("H137").Value = 3336209.5
$ws.Range("I137").Value = 5557588.5
$ws.Range("K137").Value = 16672765.5
$ws.Range("M137").Value = -16670215.5
$ws.Range("H138").Value = 2489131.8
$ws.Range("I138").Value = 905.625
$ws.Range("J138").Value = 6175392.5
$ws.Range("K138").Value = 2716.875
$ws.Range("L138").Value = 18526177.5
$ws.Range("M138").Value = 2423.125
$ws.Range("N138").Value = -18536457.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2070.93
$ws.Range("I32").Value = 1759.3478
$ws.Range("J32").Value = 2764.4517
$ws.Range("K32").Value = 1759.3478
$ws.Range("L32").Value = 2764.4517
$ws.Range("M32").Value = -1472.3478
$ws.Range("N32").Value = -3338.4517
$ws.Range("H45").Value = 1653.3846
$ws.Range("I45").Value = 985
$ws.Range("K45").Value = 985
$ws.Range("M45").Value = -608
$ws.Range("H74").Value = 7002367
$ws.Range("I74").Value = 9297933
$ws.Range("J74").Value = 115668.22
$ws.Range("K74").Value = 9297933
$ws.Range("L74").Value = 115668.22
$ws.Range("M74").Value = -9297059
$ws.Range("N74").Value = -117416.22
$ws.Range("H77").Value = 7002367
$ws.Range("I77").Value = 9297933
$ws.Range("J77").Value = 115668.22
$ws.Range("K77").Value = 46489665
$ws.Range("L77").Value = 578341.1
$ws.Range("M77").Value = -46485297
$ws.Range("N77").Value = -587077.1
$ws.Range("H122").Value = 5850269
$ws.Range("I122").Value = 2514.2856
$ws.Range("J122").Value = 22223982
$ws.Range("K122").Value = 7542.8568
$ws.Range("L122").Value = 66671946
$ws.Range("M122").Value = -5092.8568
$ws.Range("N122").Value = -66676846

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 886
$ws.Range("I64").Value = 570.8570999999999
$ws.Range("J64").Value = 1253.6666
$ws.Range("K64").Value = 570.8570999999999
$ws.Range("L64").Value = 1253.6666
$ws.Range("M64").Value = -345.8570999999999
$ws.Range("N64").Value = -1703.6666
$ws.Range("H67").Value = 886
$ws.Range("I67").Value = 570.8570999999999
$ws.Range("J67").Value = 1253.6666
$ws.Range("K67").Value = 570.8570999999999
$ws.Range("L67").Value = 1253.6666
$ws.Range("M67").Value = 209.1429000000001
$ws.Range("N67").Value = -2813.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 7295.421
$ws.Range("I99").Value = 7442.7856
$ws.Range("J99").Value = 6882.8
$ws.Range("K99").Value = 7442.7856
$ws.Range("L99").Value = 6882.8
$ws.Range("M99").Value = -5944.7856
$ws.Range("N99").Value = -9878.799999999999
$ws.Range("H126").Value = 7295.421
$ws.Range("I126").Value = 7442.7856
$ws.Range("J126").Value = 6882.8
$ws.Range("K126").Value = 22328.3568
$ws.Range("L126").Value = 20648.4
$ws.Range("M126").Value = -19858.3568
$ws.Range("N126").Value = -25588.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 367.75
$ws.Range("I40").Value = 51.625
$ws.Range("K40").Value = 206.5
$ws.Range("M40").Value = -137.5
$ws.Range("H103").Value = 3997
$ws.Range("I103").Value = 843.3333
$ws.Range("J103").Value = 4857.091
$ws.Range("K103").Value = 2529.9999
$ws.Range("L103").Value = 14571.273
$ws.Range("M103").Value = -1650.9999
$ws.Range("N103").Value = -16329.273
$ws.Range("H131").Value = 1031.7561
$ws.Range("I131").Value = 426.41666
$ws.Range("J131").Value = 1135.5286
$ws.Range("K131").Value = 1279.24998
$ws.Range("L131").Value = 3406.5858
$ws.Range("M131").Value = 3760.75002
$ws.Range("N131").Value = -13486.5858
$ws.Range("H137").Value = 49373
$ws.Range("I137").Value = 3000
$ws.Range("J137").Value = 57804.453
$ws.Range("K137").Value = 9000
$ws.Range("L137").Value = 173413.359
$ws.Range("M137").Value = -3900
$ws.Range("N137").Value = -183613.359

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5423.2666
$ws.Range("I70").Value = 5215.737
$ws.Range("K70").Value = 5215.737
$ws.Range("M70").Value = -4945.737
$ws.Range("H73").Value = 5423.2666
$ws.Range("I73").Value = 5215.737
$ws.Range("K73").Value = 5215.737
$ws.Range("M73").Value = -4279.737
$ws.Range("H122").Value = 3028.2856
$ws.Range("I122").Value = 2675.68
$ws.Range("K122").Value = 8027.039999999999
$ws.Range("M122").Value = -5577.039999999999
$ws.Range("H123").Value = 23357.8
$ws.Range("J123").Value = 23357.8
$ws.Range("L123").Value = 23357.8
$ws.Range("N123").Value = -28257.8
$ws.Range("H132").Value = 68731.42999999999
$ws.Range("I132").Value = 49193.523
$ws.Range("J132").Value = 114319.89
$ws.Range("K132").Value = 147580.569
$ws.Range("L132").Value = 342959.67
$ws.Range("M132").Value = -145050.569
$ws.Range("N132").Value = -348019.67

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3756
$ws.Range("J122").Value = 3827.5
$ws.Range("L122").Value = 11482.5
$ws.Range("N122").Value = -16382.5
$ws.Range("H132").Value = 20301.219
$ws.Range("I132").Value = 9201.617
$ws.Range("J132").Value = 74213.57000000001
$ws.Range("K132").Value = 27604.851
$ws.Range("L132").Value = 222640.71
$ws.Range("M132").Value = -25074.851
$ws.Range("N132").Value = -227700.71
$ws.Range("H136").Value = 16223.091
$ws.Range("I136").Value = 23410.8
$ws.Range("J136").Value = 10233.333
$ws.Range("K136").Value = 70232.39999999999
$ws.Range("L136").Value = 30699.999
$ws.Range("M136").Value = -67682.39999999999
$ws.Range("N136").Value = -35799.999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2400.054
$ws.Range("I122").Value = 2059.2273
$ws.Range("J122").Value = 2899.9333
$ws.Range("K122").Value = 6177.6819
$ws.Range("L122").Value = 8699.7999
$ws.Range("M122").Value = -3727.6819
$ws.Range("N122").Value = -13599.7999
$ws.Range("H136").Value = 45817.535
$ws.Range("I136").Value = 31025.242
$ws.Range("K136").Value = 93075.726
$ws.Range("M136").Value = -90525.726
